# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# The K column (column G) values are recalculated; update the affected cells
# in place with their newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    4  = 2
    5  = 0
    6  = 0
    7  = 1
    8  = 0
    9  = 1
    10 = 0
    11 = 1
    12 = 0
    13 = 2
    15 = 1
    16 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
